$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# This edit performs a 3-way rotation of certain field values among rows 3, 4 and 6:
#   new row3 <- old row6 (Id/location/accuracy fields)
#   new row4 <- old row3 (Id/location/accuracy/reporter fields, plus the "empty marker" cells)
#   new row6 <- old row4 (Id/location/accuracy/reporter fields, minus the "empty marker" cells)

# --- capture all needed "before" values first ---
$row3_A  = $ws.Range("A3").Value2
$row3_P  = $ws.Range("P3").Value2
$row3_Q  = $ws.Range("Q3").Value2
$row3_R  = $ws.Range("R3").Value2
$row3_S  = $ws.Range("S3").Value2
$row3_AW = $ws.Range("AW3").Value2
$row3_AX = $ws.Range("AX3").Value2

$row4_A  = $ws.Range("A4").Value2
$row4_P  = $ws.Range("P4").Value2
$row4_Q  = $ws.Range("Q4").Value2
$row4_R  = $ws.Range("R4").Value2
$row4_S  = $ws.Range("S4").Value2
$row4_AW = $ws.Range("AW4").Value2
$row4_AX = $ws.Range("AX4").Value2

$row6_A  = $ws.Range("A6").Value2
$row6_P  = $ws.Range("P6").Value2
$row6_Q  = $ws.Range("Q6").Value2
$row6_R  = $ws.Range("R6").Value2
$row6_S  = $ws.Range("S6").Value2
$row6_AW = $ws.Range("AW6").Value2
$row6_AX = $ws.Range("AX6").Value2

# --- add the blank "marker" cells to row 4 (mirroring row 3's layout) before row 3 changes ---
$ws.Range("J3").Copy($ws.Range("J4"))
$ws.Range("L3").Copy($ws.Range("L4"))
$ws.Range("N3").Copy($ws.Range("N4"))
$ws.Range("AF3").Copy($ws.Range("AF4"))

# --- row 3 gets row 6's values ---
$ws.Range("A3").Value2 = $row6_A
$ws.Range("P3").Value2 = $row6_P
$ws.Range("Q3").Value2 = $row6_Q
$ws.Range("R3").Value2 = $row6_R
$ws.Range("S3").Value2 = $row6_S

# --- row 4 gets row 3's values ---
$ws.Range("A4").Value2  = $row3_A
$ws.Range("P4").Value2  = $row3_P
$ws.Range("Q4").Value2  = $row3_Q
$ws.Range("R4").Value2  = $row3_R
$ws.Range("S4").Value2  = $row3_S
$ws.Range("AW4").Value2 = $row3_AW
$ws.Range("AX4").Value2 = $row3_AX

# --- row 6 gets row 4's values, and loses the blank "marker" cells (mirroring row 4's old layout) ---
$ws.Range("A6").Value2  = $row4_A
$ws.Range("P6").Value2  = $row4_P
$ws.Range("Q6").Value2  = $row4_Q
$ws.Range("R6").Value2  = $row4_R
$ws.Range("S6").Value2  = $row4_S
$ws.Range("AW6").Value2 = $row4_AW
$ws.Range("AX6").Value2 = $row4_AX

$ws.Range("J6").ClearContents()
$ws.Range("L6").ClearContents()
$ws.Range("N6").ClearContents()
$ws.Range("AF6").ClearContents()
